$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 4 hold the same credential pair; refresh it with the new
# generated manager username / password values.
$ws.Range("A2").Value = "mngr365881"
$ws.Range("B2").Value = "jYmebUz"
$ws.Range("A4").Value = "mngr365881"
$ws.Range("B4").Value = "jYmebUz"

# Leave the selection where it was saved in the workbook view state.
$ws.Range("J16").Select()
